# ENH: adding references to CV and Bibliography and fixing citation
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the blank paragraph between "Dear Colleagues," and the bold
#    "Citation:" paragraph.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "`r") {
        $para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2) Citation paragraph: "...diverse phenomena of time travel, temporal
#    paradox, and ionized hydrogen and helium of interstellar origin."
#    becomes "...diverse phenomena of space-time travel, temporal
#    paradox, and Exo-Biology."
# ---------------------------------------------------------------------
$old = "sustained impact and cross-disciplinary breakthroughs in diverse phenomena of time travel, temporal paradox, and ionized hydrogen and helium of interstellar origin."
$new = "sustained impact and cross-disciplinary breakthroughs in diverse phenomena of space-time travel, temporal paradox, and Exo-Biology."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# ---------------------------------------------------------------------
# 3) Add a citation to the Bibliography right after "...many authors."
# ---------------------------------------------------------------------
$old = "and involving dozens of research groups and many authors."
$new = "and involving dozens of research groups and many authors [see CV section Collaborations and Service activities]."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# ---------------------------------------------------------------------
# 4) Rework the opening of the diversity-service paragraph.
# ---------------------------------------------------------------------
$old = "There is one service contribution to the community that I feel is outstanding and truly deserves recognition: John" + [char]0x2019 + "s commitment to, and promotion of, diversity."
$new = "John" + [char]0x2019 + "s service contribution to the community truly deserves recognition: John" + [char]0x2019 + "s commitment to, and promotion of, diversity."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# ---------------------------------------------------------------------
# 5) Add a citation to the CV presentations after "...he deserves."
# ---------------------------------------------------------------------
$old = "he is not unduly given more credit than he deserves."
$new = "he is not unduly given more credit than he deserves [See CV section Presentations " + [char]0x2013 + " " + [char]0x201C + "How Using Team Science Ensured Safe Space-Time Travel" + [char]0x201D + " Plenary at the AGU meeting 1920, and 2020]."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# ---------------------------------------------------------------------
# 6) Add "Prof. " before the signature name, and drop the trailing
#    department/university address.
# ---------------------------------------------------------------------
$old = "D. R. Who, PhD" + [char]0xA0 + ", Time travel Department, University of Gallifrey"
$new = "Prof. D. R. Who, PhD" + [char]0xA0
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
